$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns A (type) and C (code) for rows 2..19 (row 1 = header)
$typeVals = @("migration","installation","portabilité","migration","installation","migration","migration","installation","portabilité","migration","installation","migration","portabilité","portabilité","migration","portabilité","portabilité","portabilité")
$codeVals = @("as-xxxx-xx","ot-xxxx-xx","xx-ftth-xx"," xx-ftth-xx","ot-xxxx-xx","xx-rsta-xx","as-xxxx-xx","as-xxxx-xx","xx-ftth-xx","xx-rsta-xx","ot-xxxx-xx","xx-rsta-xx","xx-ftth-xx","ot-xxxx-xx","xx-rsta-xx","as-xxxx-xx","as-xxxx-xx","xx-ftth-xx")

for ($i = 0; $i -lt 18; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $typeVals[$i]
    $ws.Cells.Item($row, 3).Value = $codeVals[$i]
    $ws.Cells.Item($row, 4).Value = 45231
    $hour = $i + 1
    $timeCell = $ws.Cells.Item($row, 5)
    $timeCell.Value = $hour / 24.0
    $timeCell.NumberFormat = "h:mm:ss"
}

$ws.Range("A3").Select()
